$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

# The 58 new customer id values (shared-string order matches the diff),
# which will populate rows 19..76 (A19:A76) in order.
$values = @(
    "cus_Kv3PKhR10oeVv8",
    "cus_Kv3PFG1hgRZWjb",
    "cus_Kv3PVe4TFfsmn1",
    "cus_Kv3PBAyoMxQq0o",
    "cus_Kv3PaJj4DnYMN5",
    "cus_Kv3PSwErcjCMlm",
    "cus_Kv3PNceLuZunTQ",
    "cus_Kv3PhkC7XR8mZJ",
    "cus_Kv3PsbsP7njhSP",
    "cus_Kv3PV0lNudyNh0",
    "cus_Kv3P6rTSOR4Rs2",
    "cus_Kv3PqPLqSMMPUL",
    "cus_Kv3PI8zu2CJVkT",
    "cus_Kv3OzZlLr69Caq",
    "cus_Kv3OLHTvvWS8fP",
    "cus_Kv3ONbCYIbBmEd",
    "cus_Kv3OWODvgmC3zn",
    "cus_Kv3OOuXiVaJMbb",
    "cus_Kv3O9eiZXUUj5X",
    "cus_Kv3OAbnzNA4lvv",
    "cus_Kv3Oc7sbKBB19j",
    "cus_Kv3ODVudIqhtUk",
    "cus_Kv3OyPLrN82gZM",
    "cus_Kv3O9wlhwerkLy",
    "cus_Kv3O8oXDmF9er1",
    "cus_Kv3OYkh6rdGviM",
    "cus_Kv3OpcSw4nRwYR",
    "cus_Kv3OEUZtAkNf1o",
    "cus_Kv3OAZXPwYPNjO",
    "cus_Kv3OJTfbQY3Fxy",
    "cus_Kv3Opqi4Uy7Uc1",
    "cus_Kv3OPZUxy6Cszz",
    "cus_Kv3O3o0IsQr3qH",
    "cus_Kv3OWXQAUqEedz",
    "cus_Kv3OBE59AKJ4UB",
    "cus_Kv3OTW7Hu3Rf3s",
    "cus_Kv3OHYq3WOKl14",
    "cus_Kv3O0O1LTfoynI",
    "cus_Kv3O6ckulx1QxD",
    "cus_Kv3OaUog44fgIK",
    "cus_Kv3OSuyNl4aZ8D",
    "cus_Kv3OAFDXmPhCjC",
    "cus_Kv3OKoiJlFNvyH",
    "cus_Kv3NN05uR9X9zN",
    "cus_Kv3NEGHcmOSdid",
    "cus_Kv3NPF09fFXY32",
    "cus_Kv3N2Ukbe1AXax",
    "cus_Kv3NoPGPPz2Jlm",
    "cus_Kv3NDr81Lqgiej",
    "cus_Kv3NfYE6IOX4Pe",
    "cus_Kv3NNARjxmimdQ",
    "cus_Kv3Nk81SMdSLDL",
    "cus_Kv3NDSOTqzfNhd",
    "cus_Kv3Nwu4r8v76xV",
    "cus_Kv3Nib3fHFBabm",
    "cus_Kv3NLltNwlxwlr",
    "cus_Kv3NpNnfkVwng2",
    "cus_Kv3NwW2j0BdYwY"
)

$startRow = 19
$lastExistingRow = 68

# Rows 19..68 already exist and already carry the bordered style (s="1"),
# so a plain value write is enough - it keeps reusing that existing style.
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    if ($row -le $lastExistingRow) {
        $ws.Cells.Item($row, 1).Value = $values[$i]
    }
}

# Rows 69..76 are brand new, so first write their values, then clone the
# formatting (thin border) from the existing A19 cell via copy/paste-special
# so we reuse the existing style entry instead of minting a new one.
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    if ($row -gt $lastExistingRow) {
        $ws.Cells.Item($row, 1).Value = $values[$i]
    }
}

$ws.Range("A19").Copy()
$ws.Range("A69:A76").PasteSpecial(-4122)

# Update the view: scroll position + selection to mirror the authored diff.
$ws.Range("A49:A76").Select()
$excel.ActiveWindow.ScrollRow = 37
